$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial (45171 = 2023-09-02) for every
# data row (rows 2 through 416). Bump it by one day (45172 = 2023-09-03)
# for all of them.
$ws.Range("C2:C416").Value = 45172
